# Refresh the cryptos list (prices / 1h volume %, and a couple of
# ranking swaps) to match the latest GitHub Actions data pull.
#
# Note: several "Price" cells (column D) look like plain numbers
# (e.g. "604.68"), and Excel would otherwise auto-convert them from
# text to a floating point number on assignment (losing the exact
# formatting, e.g. "7.50" -> 7.5). To keep them as text exactly as
# scraped, we temporarily force the cell to Text format, assign the
# value, then restore the style back to Normal so no stray style is
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.912.12'
$ws.Range('E2').Value = '  +1.60%  '

$ws.Range('D3').Value = '3.214.76'
$ws.Range('E3').Value = '  +1.54%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '604.68'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.39%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.76'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.23%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').Value = '3.213.60'
$ws.Range('E8').Value = '  +1.47%  '

$ws.Range('E9').Value = '  +0.37%  '

$ws.Range('E10').Value = '  -0.87%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.16'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.06%  '

$ws.Range('E12').Value = '  +1.85%  '

$ws.Range('E13').Value = '  +1.86%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '39.07'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.13%  '

$ws.Range('D15').Value = '3.740.20'
$ws.Range('E15').Value = '  +1.43%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.50'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +4.88%  '

$ws.Range('D17').Value = '66.067.37'
$ws.Range('E17').Value = '  +1.67%  '

$ws.Range('D18').Value = '3.211.53'
$ws.Range('E18').Value = '  +1.40%  '

$ws.Range('E19').Value = '  +0.04%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '509.52'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.19%  '

$ws.Range('E21').Value = '  +4.65%  '

$ws.Range('E22').Value = '  +1.38%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '15.37'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.33%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.09'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.95%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '85.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.06%  '

$ws.Range('E26').Value = '  -0.09%  '

$ws.Range('E27').Value = '  +3.30%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.18'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.09%  '

$ws.Range('E29').Value = '  +4.18%  '

$ws.Range('E30').Value = '  +3.18%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '28.15'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.33%  '

$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.78'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +7.71%  '

$ws.Range('E33').Value = '  +1.76%  '

$ws.Range('E34').Value = '  -0.05%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.60'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.61%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '55.04'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.29%  '

$ws.Range('E37').Value = '  +0.98%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '480.21'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.50%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0420'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.06%  '

$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.92'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.42%  '

$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.94'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -5.14%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.298'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.02%  '

$ws.Range('E43').Value = '  +1.06%  '

$ws.Range('D44').Value = '2.946.90'
$ws.Range('E44').Value = '  -3.68%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.45'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.43%  '

$ws.Range('D46').Value = '0.0₃0640'
$ws.Range('E46').Value = '  +6.61%  '

$ws.Range('E47').Value = '  -2.04%  '

$ws.Range('E49').Value = '  +1.32%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.30'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.73%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '121.50'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.41%  '
